# Atualização da Ata de Reuniões
# Adds the "09/10 (Terça-Feira)" meeting row (13) and the new
# "10/10 (Quarta-Feira)" header row (14), bumps the first meeting date
# in row 12 from 02/10 to 08/10, resizes a few rows, and re-selects B19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Row 12: the sprint's first day moves from 02/10 to 08/10
# ------------------------------------------------------------------
$ws.Range("A12").Value2 = "08/10 (Segunda-Feira)"
$ws.Rows.Item(12).RowHeight = 144.75

# ------------------------------------------------------------------
# 2. Row 13: brand-new meeting entry (09/10 - Terça-Feira)
#    Start with the same look as row 12 (date cell, two time cells,
#    "members present" cell) and then build the two new styles the
#    new cells in F13/G13 need.
# ------------------------------------------------------------------
$ws.Range("A9:A9").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").NumberFormat = "d-mmm"
$ws.Range("A13").Value2 = "09/10 (Terça-Feira)"

$ws.Range("B9").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$ws.Range("B13").Value2 = 0.70138888888888884
$ws.Range("C13").Value2 = 0.71180555555555547

$ws.Range("D8").Copy()
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("D13").Value2 = "Arthur`nDaniel`nDavi`nLeonardo`nMatteus`nPedro"

$ws.Range("D12").Copy()
$ws.Range("F13").PasteSpecial(-4122)
$ws.Range("F13").HorizontalAlignment = 1
$ws.Range("F13").Value2 = "foi reforçados com os  membros ausentes na ultima reunião os conteudos que serão desenvolvidos na semana e o que cada membro vai desenvolver; "

$ws.Range("E12").Copy()
$ws.Range("G13").PasteSpecial(-4122)
$ws.Range("G13").Value2 = "Validamos e reforçamos os prazos de entrega da semana ;"

$ws.Rows.Item(13).RowHeight = 138

# ------------------------------------------------------------------
# 3. Row 14: new trailing date header (10/10 - Quarta-Feira)
# ------------------------------------------------------------------
$ws.Range("A9:C9").Copy()
$ws.Range("A14:C14").PasteSpecial(-4122)
$ws.Range("A14").Value2 = "10/10 (Quarta-Feira)"
$ws.Range("B14").Value2 = ""
$ws.Range("C14").Value2 = ""

$ws.Rows.Item(14).RowHeight = 112.5

# ------------------------------------------------------------------
# 4. Minor row-height tweaks on the existing entries above
# ------------------------------------------------------------------
$ws.Rows.Item(7).RowHeight = 126.75
$ws.Rows.Item(8).RowHeight = 144.75
$ws.Rows.Item(9).RowHeight = 108.75

# ------------------------------------------------------------------
# 5. Selection moved from F19 to B19
# ------------------------------------------------------------------
$ws.Range("B19").Select()
